# Update "想去人数" (want-to-go count) values in column F for the
# "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet.
#
# 展览 sheet: rows 3-6 correspond to F3,F4,F5,F6
# 全部类型 sheet: same events, but row numbering differs (F3,F4,F5,F7)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 81
$ws1.Range("F4").Value = 2197
$ws1.Range("F5").Value = 193
$ws1.Range("F6").Value = 366

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 81
$ws4.Range("F4").Value = 2197
$ws4.Range("F5").Value = 193
$ws4.Range("F7").Value = 366
